$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.482.91"
$ws.Range("E2").Value = "  +2.37%  "

$ws.Range("D3").Value = "3.189.84"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.559"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.13%  "

$ws.Range("D9").Value = "3.187.64"
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("E10").Value = "  +1.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.91"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.44%  "

$ws.Range("E12").Value = "  +3.54%  "

$ws.Range("E13").Value = "  +1.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.56%  "

$ws.Range("D15").Value = "3.711.80"
$ws.Range("E15").Value = "  +1.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.12%  "

$ws.Range("D17").Value = "66.509.24"
$ws.Range("E17").Value = "  +2.49%  "

$ws.Range("D18").Value = "3.189.94"
$ws.Range("E18").Value = "  +1.01%  "

$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "521.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.741"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.85%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.37%  "

$ws.Range("E28").Value = "  +2.99%  "

$ws.Range("E29").Value = "  +7.33%  "

$ws.Range("E30").Value = "  +13.58%  "

$ws.Range("E31").Value = "  +5.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.66%  "

$ws.Range("E33").Value = "  +2.75%  "

$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "511.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "54.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0905"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.54%  "

$ws.Range("E39").Value = "  +2.19%  "

$ws.Range("E40").Value = "  +11.09%  "

$ws.Range("E41").Value = "  +1.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.99%  "

$ws.Range("E43").Value = "  +6.75%  "

$ws.Range("D44").Value = "0.0₃0675"
$ws.Range("E44").Value = "  +14.87%  "

$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("D46").Value = "2.901.84"
$ws.Range("E46").Value = "  -3.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.03%  "

$ws.Range("E48").Value = "  +11.31%  "

$ws.Range("E49").Value = "  +3.59%  "

$ws.Range("E51").Value = "  +5.23%  "

